$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-14 Thursday" "2024-03-15 Friday"

Replace-Text "89×80=7120" "19×52=988"
Replace-Text "78×52=4056" "30×94=2820"
Replace-Text "21×53=1113" "53×50=2650"
Replace-Text "50×53=2650" "85×96=8160"
Replace-Text "85×65=5525" "18×71=1278"

Replace-Text "48×77=3696" "31×48=1488"
Replace-Text "81×19=1539" "30×40=1200"
Replace-Text "69×22=1518" "21×97=2037"
Replace-Text "15×57=855" "47×88=4136"
Replace-Text "49×30=1470" "49×81=3969"

Replace-Text "76×22=1672" "51×39=1989"
Replace-Text "47×91=4277" "99×27=2673"
Replace-Text "25×66=1650" "89×59=5251"
Replace-Text "59×50=2950" "13×80=1040"
Replace-Text "71×86=6106" "25×38=950"

Replace-Text "52×46=2392" "95×81=7695"
Replace-Text "64×54=3456" "39×73=2847"
Replace-Text "37×26=962" "63×75=4725"
Replace-Text "30×46=1380" "29×30=870"
Replace-Text "11×22=242" "79×50=3950"

Replace-Text "82×51=4182" "40×35=1400"
Replace-Text "49×68=3332" "47×64=3008"
Replace-Text "61×60=3660" "82×99=8118"
Replace-Text "65×90=5850" "52×59=3068"
Replace-Text "91×59=5369" "24×26=624"
